$wb = $excel.ActiveWorkbook

$wsTests = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

# Add the new workflow entry row to the "Tests" sheet
$wsTests.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsTests.Range("B10").Value = "Success"

# Add the same new workflow entry row to the "Result" sheet
$wsResult.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsResult.Range("B10").Value = "Success"

# Update selections / active sheet to match target state
$wsResult.Range("A10").Select()
$wsTests.Activate()
$wsTests.Range("B20").Select()
